# Updated cryptos list on Sun Mar 31 03:25:52 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto ranking sheet, and re-ranks two pairs of coins whose relative
# order flipped (WrappedEther/TRON at rows 19-20, ThetaToken/ApeXProtocol
# at rows 46-47) by swapping their Coin/Link/Price/Volume cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is one cell update. Price cells (column D) are flagged
# ForceText so numeric-looking strings (e.g. "1.00", "53.20", "4.20")
# keep their exact textual form -- using NumberFormat "@" while writing
# the value and then ClearFormats() afterwards prevents Excel's automatic
# "looks like a number" coercion from dropping trailing zeros / collapsing
# multi-dot price strings like "69.687.02", while leaving the cell's
# style/format untouched once the write is done.
$updates = @(
    @{ Cell = "D2"; Value = "69.687.02"; ForceText = $true },
    @{ Cell = "E2"; Value = "  -0.14%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "3.524.98"; ForceText = $true },
    @{ Cell = "E3"; Value = "  +0.92%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  -0.07%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "606.03"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +0.01%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "193.87"; ForceText = $true },
    @{ Cell = "E6"; Value = "  +0.95%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "0.623"; ForceText = $true },
    @{ Cell = "E7"; Value = "  -0.52%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E8"; Value = "  +0.04%  "; ForceText = $false },
    @{ Cell = "D9"; Value = "0.201"; ForceText = $true },
    @{ Cell = "E9"; Value = "  -6.18%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.645"; ForceText = $true },
    @{ Cell = "E10"; Value = "  -2.31%  "; ForceText = $false },
    @{ Cell = "D11"; Value = "53.20"; ForceText = $true },
    @{ Cell = "E11"; Value = "  -0.44%  "; ForceText = $false },
    @{ Cell = "E12"; Value = "  -1.82%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "9.45"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -1.51%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "4.087.96"; ForceText = $true },
    @{ Cell = "E14"; Value = "  +0.56%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "592.07"; ForceText = $true },
    @{ Cell = "E15"; Value = "  -4.20%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "69.789.97"; ForceText = $true },
    @{ Cell = "E16"; Value = "  -0.22%  "; ForceText = $false },
    @{ Cell = "E17"; Value = "  +0.20%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "18.95"; ForceText = $true },
    @{ Cell = "E18"; Value = "  +0.61%  "; ForceText = $false },
    @{ Cell = "B19"; Value = "TRON"; ForceText = $false },
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; ForceText = $false },
    @{ Cell = "D19"; Value = "0.122"; ForceText = $true },
    @{ Cell = "E19"; Value = "  +1.78%  "; ForceText = $false },
    @{ Cell = "B20"; Value = "WrappedEther"; ForceText = $false },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; ForceText = $false },
    @{ Cell = "D20"; Value = "3.507.05"; ForceText = $true },
    @{ Cell = "E20"; Value = "  +0.09%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "0.981"; ForceText = $true },
    @{ Cell = "E21"; Value = "  -0.87%  "; ForceText = $false },
    @{ Cell = "D22"; Value = "17.73"; ForceText = $true },
    @{ Cell = "E22"; Value = "  -0.72%  "; ForceText = $false },
    @{ Cell = "E23"; Value = "  +1.62%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "102.46"; ForceText = $true },
    @{ Cell = "E24"; Value = "  -2.77%  "; ForceText = $false },
    @{ Cell = "D25"; Value = "4.62"; ForceText = $true },
    @{ Cell = "E25"; Value = "  -0.37%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "3.03"; ForceText = $true },
    @{ Cell = "E26"; Value = "  -0.22%  "; ForceText = $false },
    @{ Cell = "D27"; Value = "10.73"; ForceText = $true },
    @{ Cell = "E27"; Value = "  -2.34%  "; ForceText = $false },
    @{ Cell = "D28"; Value = "9.49"; ForceText = $true },
    @{ Cell = "E28"; Value = "  -3.71%  "; ForceText = $false },
    @{ Cell = "D29"; Value = "33.09"; ForceText = $true },
    @{ Cell = "E29"; Value = "  -3.49%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "7.01"; ForceText = $true },
    @{ Cell = "E30"; Value = "  -1.53%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "4.20"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -1.78%  "; ForceText = $false },
    @{ Cell = "D32"; Value = "12.30"; ForceText = $true },
    @{ Cell = "E32"; Value = "  -2.27%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "0.114"; ForceText = $true },
    @{ Cell = "E33"; Value = "  -0.28%  "; ForceText = $false },
    @{ Cell = "D34"; Value = "63.27"; ForceText = $true },
    @{ Cell = "E34"; Value = "  -1.27%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "3.793.18"; ForceText = $true },
    @{ Cell = "E35"; Value = "  +1.46%  "; ForceText = $false },
    @{ Cell = "E36"; Value = "  +3.15%  "; ForceText = $false },
    @{ Cell = "E37"; Value = "  +0.22%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "0.0₃0806"; ForceText = $true },
    @{ Cell = "E38"; Value = "  +1.61%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "511.40"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -2.50%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "0.390"; ForceText = $true },
    @{ Cell = "E40"; Value = "  +0.05%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "3.57"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -0.25%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "36.46"; ForceText = $true },
    @{ Cell = "E42"; Value = "  -0.61%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "0.133"; ForceText = $true },
    @{ Cell = "E43"; Value = "  -2.91%  "; ForceText = $false },
    @{ Cell = "D44"; Value = "0.0446"; ForceText = $true },
    @{ Cell = "E44"; Value = "  -3.59%  "; ForceText = $false },
    @{ Cell = "E45"; Value = "  -0.92%  "; ForceText = $false },
    @{ Cell = "B46"; Value = "ApeXProtocol"; ForceText = $false },
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"; ForceText = $false },
    @{ Cell = "D46"; Value = "3.31"; ForceText = $true },
    @{ Cell = "E46"; Value = "  -0.07%  "; ForceText = $false },
    @{ Cell = "B47"; Value = "ThetaToken"; ForceText = $false },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"; ForceText = $false },
    @{ Cell = "D47"; Value = "2.80"; ForceText = $true },
    @{ Cell = "E47"; Value = "  -2.16%  "; ForceText = $false },
    @{ Cell = "E48"; Value = "  +0.02%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "8.45"; ForceText = $true },
    @{ Cell = "E49"; Value = "  -3.35%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "0.000245"; ForceText = $true },
    @{ Cell = "E50"; Value = "  +3.32%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "1.31"; ForceText = $true },
    @{ Cell = "E51"; Value = "  +1.90%  "; ForceText = $false }

)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $u.Value
    }
}
